# Updated cryptos list on Fri Oct 11 14:59:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "61.820.97"
Set-TextValue "E2" "  +2.03%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.413.45"
Set-TextValue "E3" "  +0.20%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.57%  "

# Row 5 - BNB
Set-TextValue "D5" "569.45"
Set-TextValue "E5" "  +0.92%  "

# Row 6 - Solana
Set-TextValue "D6" "144.09"
Set-TextValue "E6" "  +4.77%  "

# Row 7 - USDC
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.49%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.426.65"

# Row 10 - Dogecoin
Set-TextValue "E10" "  +3.82%  "

# Row 11 - TRON
Set-TextValue "E11" "  +0.70%  "

# Row 12 - Toncoin
Set-TextValue "E12" "  +3.88%  "

# Row 13 - Cardano
Set-TextValue "E13" "  +3.79%  "

# Row 14 - Avalanche
Set-TextValue "D14" "26.61"
Set-TextValue "E14" "  +3.95%  "

# Row 15 - ShibaInu
Set-TextValue "E15" "  +5.11%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "2.835.56"
Set-TextValue "E16" "  +0.44%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "61.606.48"
Set-TextValue "E17" "  +1.57%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.404.00"
Set-TextValue "E18" "  +0.11%  "

# Row 19 - Uniswap
Set-TextValue "D19" "7.97"
Set-TextValue "E19" "  -0.57%  "

# Row 20 - Chainlink
Set-TextValue "D20" "10.73"
Set-TextValue "E20" "  +2.13%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "325.42"
Set-TextValue "E21" "  +0.98%  "

# Row 22 - Polkadot
Set-TextValue "E22" "  +1.52%  "

# Row 23 - was SuiNetwork, becomes LEO
Set-TextValue "B23" "LEO"
Set-TextValue "C23" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "6.07"
Set-TextValue "E23" "  -1.77%  "

# Row 24 - was LEO, becomes SuiNetwork
Set-TextValue "B24" "SuiNetwork"
Set-TextValue "C24" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D24" "2.03"
Set-TextValue "E24" "  +12.39%  "

# Row 25 - Dai
Set-TextValue "E25" "  -0.17%  "

# Row 26 - Litecoin
Set-TextValue "D26" "65.24"
Set-TextValue "E26" "  +1.80%  "

# Row 27 - Bittensor
Set-TextValue "D27" "620.28"
Set-TextValue "E27" "  +11.80%  "

# Row 28 - Aptos
Set-TextValue "D28" "8.40"
Set-TextValue "E28" "  +2.17%  "

# Row 29 - PEPE
Set-TextValue "D29" "0.0₃0967"
Set-TextValue "E29" "  +5.36%  "

# Row 30 - WrappedeETH
Set-TextValue "E30" "  -0.21%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "8.03"
Set-TextValue "E31" "  +1.91%  "

# Row 32 - Fetch.AI
Set-TextValue "E32" "  +8.81%  "

# Row 33 - PancakeSwap
Set-TextValue "E33" "  +1.01%  "

# Row 34 - Kaspa
Set-TextValue "E34" "  +2.24%  "

# Row 35 - ImmutableX
Set-TextValue "E35" "  +5.06%  "

# Row 36 - FirstDigitalUSD
Set-TextValue "D36" "0.997"
Set-TextValue "E36" "  -0.60%  "

# Row 37 - Monero
Set-TextValue "D37" "153.12"
Set-TextValue "E37" "  +0.90%  "

# Row 38 - NEARProtocol
Set-TextValue "E38" "  +2.42%  "

# Row 40 - RenderToken
Set-TextValue "D40" "5.36"
Set-TextValue "E40" "  +5.88%  "

# Row 41 - EthereumClassic
Set-TextValue "D41" "18.43"
Set-TextValue "E41" "  +1.74%  "

# Row 42 - dogwifhat
Set-TextValue "D42" "2.62"
Set-TextValue "E42" "  +11.90%  "

# Row 43 - Stacks
Set-TextValue "E43" "  +4.50%  "

# Row 44 - was OKB, becomes USDe
Set-TextValue "B44" "USDe"
Set-TextValue "C44" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D44" "0.999"
Set-TextValue "E44" "  -0.08%  "

# Row 45 - was USDe, becomes OKB
Set-TextValue "B45" "OKB"
Set-TextValue "C45" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D45" "42.18"
Set-TextValue "E45" "  +2.79%  "

# Row 46 - BabyDogeCoin
Set-TextValue "E46" "  -2.97%  "

# Row 47 - Aave
Set-TextValue "D47" "142.52"
Set-TextValue "E47" "  +0.33%  "

# Row 48 - Filecoin
Set-TextValue "E48" "  +1.62%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "20.05"
Set-TextValue "E49" "  +5.29%  "

# Row 50 - Mantle
Set-TextValue "D50" "0.596"
Set-TextValue "E50" "  +2.59%  "

# Row 51 - Hedera
Set-TextValue "D51" "0.0513"
Set-TextValue "E51" "  +2.99%  "
